$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Marking" row correct-count (B11) and "Total" row (B12 / E12)
$ws.Range("B11").Value = 5
$ws.Range("B12").Value = 95
$ws.Range("E12").Value = "95/140"
